$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- New row 45: "Left Rotate the Array by One And Kth" ----

# A45: date serial, reuse A43's date style (numFmtId 14, centered)
$ws.Range("A43").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A45").Value = 45427

# B45: day-of-week text, same style as B43 (no explicit style)
$ws.Range("B43").Copy() | Out-Null
$ws.Range("B45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B45").Value = "WED"

# D45: difficulty text, same style as D43 (no explicit style)
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D45").Value = "Easy"

# F45: running count
$ws.Range("F43").Copy() | Out-Null
$ws.Range("F45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("F45").Value = 25

# E45: hyperlink to the GFG problem (URL text becomes the new shared string
# right before the question-text string, matching the source ordering)
$link = "https://www.geeksforgeeks.org/problems/rotate-array-by-n-elements-1587115621/1?itm_source=geeksforgeeks&itm_medium=article&itm_campaign=bottom_sticky_on_article"
$ws.Hyperlinks.Add($ws.Range("E45"), $link, "", "", $link) | Out-Null
$ws.Range("E43").Copy() | Out-Null
$ws.Range("E45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# C45: question text, same style as C43 (no explicit style)
$ws.Range("C43").Copy() | Out-Null
$ws.Range("C45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C45").Value = "Left Rotate the Array by One And Kth "

# Row height to match the rest of the wrapped-text block
$ws.Rows.Item(45).RowHeight = $ws.Rows.Item(43).RowHeight

$ws.Range("C48").Select() | Out-Null
